$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "Prix Spot": add a new column Z (09-jul) mirroring column Y's
# header style, with header text + 24 hourly price values.
# ----------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the style (bold/border/centered) from the Y1 header cell onto Z1.
$wsPrix.Range("Y1").Copy() | Out-Null
$wsPrix.Range("Z1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$wsPrix.Range("Z1").Value = "09-jul"

$prixValues = @(
    71.98999999999999,
    39.37,
    33.59,
    23.37,
    15.28,
    54.77,
    15,
    61.47,
    69.13,
    61.61,
    20.99,
    0,
    11.68,
    1.38,
    0,
    5,
    2.93,
    30.02,
    64.28,
    103.96,
    111.23,
    95.16,
    115.91,
    97.37
)

for ($i = 0; $i -lt $prixValues.Count; $i++) {
    $row = 2 + $i
    $wsPrix.Cells.Item($row, 26).Value = $prixValues[$i]
}

# ----------------------------------------------------------------------
# Sheet "Gaz": append row 23 with date 2025-07-07 and its price.
# The date must stay a plain text value (matching the existing A2:A22
# cells), not get auto-converted into a date serial number, so we
# temporarily force a text format, assign the value, then restore the
# default "Normal" style so no stray style index is left on the cell.
# ----------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A23").NumberFormat = "@"
$wsGaz.Range("A23").Value = "2025-07-07"
$wsGaz.Range("A23").Style = "Normal"
$wsGaz.Range("B23").Value = 33.4

# ----------------------------------------------------------------------
# Sheet "CO2": append row 23 with date 2025-07-07 and its price.
# ----------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A23").NumberFormat = "@"
$wsCO2.Range("A23").Value = "2025-07-07"
$wsCO2.Range("A23").Style = "Normal"
$wsCO2.Range("B23").Value = 69.95999999999999
